$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.459.50"
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = "'1.817.48"
$ws.Range('E3').Value = '  +4.26%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = "'344.12"
$ws.Range('E5').Value = '  +2.96%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('D7').Value = "'0.3842"
$ws.Range('E7').Value = '  +2.40%  '
$ws.Range('D8').Value = "'0.3561"
$ws.Range('D9').Value = "'49.00"
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').Value = "'1.241"
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('D11').Value = "'0.07809"
$ws.Range('E11').Value = '  +3.69%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = "'22.53"
$ws.Range('E13').Value = '  +9.55%  '
$ws.Range('D14').Value = "'6.641"
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').Value = "'1.814.58"
$ws.Range('E15').Value = '  +4.41%  '
$ws.Range('D16').Value = "'7.243"
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = "'0.00001131"
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D18').Value = "'0.06742"
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = "'86.97"
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').Value = "'1.001"
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = "'17.73"
$ws.Range('E21').Value = '  +5.33%  '
$ws.Range('D22').Value = "'6.613"
$ws.Range('E22').Value = '  +6.54%  '
$ws.Range('D23').Value = "'13.23"
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').Value = "'27.457.23"
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('D25').Value = "'2.473"
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').Value = "'2.708"
$ws.Range('E26').Value = '  +7.02%  '
$ws.Range('D27').Value = "'22.28"
$ws.Range('E27').Value = '  +13.10%  '
$ws.Range('D28').Value = "'1.474"
$ws.Range('E28').Value = '  +3.88%  '
$ws.Range('D29').Value = "'154.35"
$ws.Range('E29').Value = '  +0.48%  '
$ws.Range('D30').Value = "'2.019.25"
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('D31').Value = "'136.48"
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('D32').Value = "'6.425"
$ws.Range('E32').Value = '  +2.93%  '
$ws.Range('D33').Value = "'4.084"
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').Value = "'14.03"
$ws.Range('E34').Value = '  +6.14%  '
$ws.Range('D35').Value = "'0.08833"
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').Value = "'1.691"
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('D37').Value = "'5.674"
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('D38').Value = "'0.7068"
$ws.Range('E38').Value = '  +12.03%  '
$ws.Range('D39').Value = "'0.06551"
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').Value = "'0.2267"
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'9.042"
$ws.Range('E41').Value = '  +3.76%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'0.02406"
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').Value = "'1.301"
$ws.Range('E43').Value = '  +4.42%  '
$ws.Range('D44').Value = "'14.85"
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').Value = "'0.6653"
$ws.Range('E45').Value = '  +8.54%  '
$ws.Range('D46').Value = "'1.001"
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = "'3.967"
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').Value = "'2.214"
$ws.Range('E48').Value = '  +6.41%  '
$ws.Range('D49').Value = "'132.69"
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('D50').Value = "'0.07326"
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = "'81.29"
$ws.Range('E51').Value = '  +4.12%  '
